# Generate Report for Archive
#
# Refresh the localization status report: items previously shown as
# "Ready for handoff" are now "In Translation". This touches the Overview
# sheet's per-locale status columns (zh-cn / de-de) as well as the Status
# column on each per-locale detail sheet. Once the text is shorter, the
# affected status columns are resized to fit their new contents.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target column width (in Excel "characters" units) that the refreshed,
# narrower status columns should end up with once resized to fit
# "In Translation" instead of the previous, longer "Ready for handoff".
$statusColumnWidth = 12.5

# Update the status text everywhere it appears and remember which
# worksheet columns actually contained it, so only those columns get
# resized afterward.
$touchedColumns = @{}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($oldStatus -eq $cell.Text) {
            $cell.Value = $newStatus
            $touchedColumns[$ws.Name + "|" + $cell.Column] = $ws
        }
    }
}

foreach ($key in $touchedColumns.Keys) {
    $ws = $touchedColumns[$key]
    $colIndex = [int]($key.Split("|")[1])
    $ws.Columns.Item($colIndex).ColumnWidth = $statusColumnWidth
}

$wb.Save()
